# Add files via upload
# ReportsDetail.xlsx: expand Sheet1 with the full System x Invoice Status
# combination table (FMS/AFS x PARKED/POSTED), and update the selections
# on both sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet2: selection becomes the whole A1:B4 block (no single active cell) ----
# (done first so that Sheet1 ends up as the active/selected tab, matching
#  the original workbook where Sheet1 was tabSelected)
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A1:B4").Select()

# ---- Sheet1: add rows 2-5 (System / Invoice Status combinations) ----
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("A2").Value = "FMS"
$ws1.Range("B2").Value = "PARKED"

$ws1.Range("A3").Value = "FMS"
$ws1.Range("B3").Value = "POSTED"

$ws1.Range("A4").Value = "AFS"
$ws1.Range("B4").Value = "PARKED"

$ws1.Range("A5").Value = "AFS"
$ws1.Range("B5").Value = "POSTED"

# Selection moves to F13 on Sheet1 (also re-activates Sheet1 as the tab)
$ws1.Range("F13").Select()
